$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Class Statistics summary box (K/L columns) ---
$ws.Range("L6").Value = 159
$ws.Range("L7").Value = 3

# --- Swap "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System" ---
$swapCells = @("G8","G9","G10","G34","G35","G36","G60","G61","G62","G86","G87","G88","G112","G113","G114","G138","G139","G140","G164","G167","G191","G194","G218","G221","G245","G248","G272","G275","G299","G302")
foreach ($cell in $swapCells) {
    $ws.Range($cell).Value = "dnasr281@gmail.com, System"
}

# --- Per-group breakdown table (rows 21-26): Recorded/Missing counts ---
$ws.Range("O21").Value = 13
$ws.Range("P21").Value = 0

$ws.Range("O22").Value = 13
$ws.Range("P22").Value = 0

$ws.Range("O23").Value = 13
$ws.Range("P23").Value = 0

$ws.Range("O24").Value = 12
$ws.Range("P24").Value = 1

$ws.Range("O25").Value = 13
$ws.Range("P25").Value = 0

$ws.Range("O26").Value = 13
$ws.Range("P26").Value = 0

# --- Percentage text cells: writing a "NN.N%" string directly gets auto-parsed into a
# real percentage NUMBER (and a brand-new number-format style) by this engine's smart
# value coercion, which would silently change the cell's style id. To keep the cell a
# literal text value (matching the source XML's t="inlineStr"/shared-string cells) AND
# keep the original shared style (s="5") untouched, force a Text format right before
# assigning, then re-stamp the original style with a formats-only paste from an
# untouched donor cell that already carries that exact style.
$pctCells = @{
    "L9"  = "50.0%";
    "L10" = "73.5%";
    "R21" = "48.1%"; "S21" = "76.9%";
    "R22" = "48.1%"; "S22" = "76.2%";
    "R23" = "48.1%"; "S23" = "80.6%";
    "R24" = "44.4%"; "S24" = "70.5%";
    "R25" = "48.1%"; "S25" = "70.1%";
    "R26" = "48.1%"; "S26" = "61.8%";
}
foreach ($cell in $pctCells.Keys) {
    $ws.Range($cell).NumberFormat = "@"
    $ws.Range($cell).Value = $pctCells[$cell]
}
$ws.Range("L15").Copy()
foreach ($cell in $pctCells.Keys) {
    $ws.Range($cell).PasteSpecial(-4122)
}

# --- Sessions newly recorded (rows 170, 197, 224, 251, 278, 305): ---
# Copy the formatting (fill/style) from the row immediately above (already "Recorded" style)
# then overwrite the Recorded-By / Students / Status values for the now-recorded session.
$sessionRows = @(
    @{Row=170; Donor=168; RecordedBy="dnasr281@gmail.com"; Students="16/23"},
    @{Row=197; Donor=196; RecordedBy="dnasr281@gmail.com"; Students="25/30"},
    @{Row=224; Donor=223; RecordedBy="dnasr281@gmail.com"; Students="21/25"},
    @{Row=251; Donor=250; RecordedBy="dnasr281@gmail.com"; Students="21/28"},
    @{Row=278; Donor=277; RecordedBy="dnasr281@gmail.com"; Students="21/26"},
    @{Row=305; Donor=304; RecordedBy="dnasr281@gmail.com"; Students="22/29"}
)

foreach ($s in $sessionRows) {
    $ws.Range("A$($s.Donor):I$($s.Donor)").Copy()
    $ws.Range("A$($s.Row):I$($s.Row)").PasteSpecial(-4122)
    $ws.Range("G$($s.Row)").Value = $s.RecordedBy
    $ws.Range("H$($s.Row)").Value = $s.Students
    $ws.Range("I$($s.Row)").Value = "Recorded"
}
